$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value2 = '66.106.69'
$ws.Range('E2').Value2 = '  +0.30%  '
$ws.Range('D3').Value2 = '3.558.69'
$ws.Range('E3').Value2 = '  +3.34%  '
$ws.Range('E4').Value2 = '  +0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value2 = '605.76'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value2 = '  +1.27%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value2 = '144.32'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value2 = '  +1.18%  '
$ws.Range('D7').Value2 = '3.557.20'
$ws.Range('E7').Value2 = '  +3.21%  '
$ws.Range('E8').Value2 = '  +0.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value2 = '0.491'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value2 = '  +4.07%  '
$ws.Range('E10').Value2 = '  +1.42%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value2 = '7.91'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value2 = '  -2.12%  '
$ws.Range('E12').Value2 = '  +1.12%  '
$ws.Range('D13').Value2 = '4.162.15'
$ws.Range('E13').Value2 = '  +3.47%  '
$ws.Range('E14').Value2 = '  +3.05%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value2 = '29.95'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value2 = '  +0.56%  '
$ws.Range('D16').Value2 = '3.561.40'
$ws.Range('E16').Value2 = '  +3.78%  '
$ws.Range('D17').Value2 = '66.202.20'
$ws.Range('E17').Value2 = '  +0.46%  '
$ws.Range('E18').Value2 = '  -0.54%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value2 = '11.29'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value2 = '  +9.45%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value2 = '6.19'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value2 = '  +1.08%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value2 = '14.83'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value2 = '  +1.31%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value2 = '428.86'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value2 = '  +2.87%  '
$ws.Range('E23').Value2 = '  +5.40%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value2 = '79.08'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value2 = '  +2.31%  '
$ws.Range('D25').Value2 = '3.701.62'
$ws.Range('E25').Value2 = '  +3.68%  '
$ws.Range('E26').Value2 = '  +0.03%  '
$ws.Range('E27').Value2 = '  +5.22%  '
$ws.Range('E28').Value2 = '  +2.67%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value2 = '7.93'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value2 = '  -0.23%  '
$ws.Range('E30').Value2 = '  -2.64%  '
$ws.Range('E31').Value2 = '  +0.12%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value2 = '25.55'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value2 = '  +3.19%  '
$ws.Range('E33').Value2 = '  +0.00%  '
$ws.Range('D34').Value2 = '3.554.61'
$ws.Range('E34').Value2 = '  +3.43%  '
$ws.Range('E35').Value2 = '  -5.41%  '
$ws.Range('E36').Value2 = '  +0.06%  '
$ws.Range('E37').Value2 = '  +2.36%  '
$ws.Range('E38').Value2 = '  +3.87%  '
$ws.Range('E39').Value2 = '  +1.34%  '
$ws.Range('E40').Value2 = '  +0.13%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value2 = '175.93'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value2 = '  +4.08%  '
$ws.Range('E42').Value2 = '  -1.54%  '
$ws.Range('E43').Value2 = '  +2.51%  '
$ws.Range('E44').Value2 = '  +0.97%  '
$ws.Range('E45').Value2 = '  +2.01%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value2 = '46.01'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value2 = '  +1.30%  '
$ws.Range('E47').Value2 = '  +0.12%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value2 = '25.78'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value2 = '  -2.05%  '
$ws.Range('E49').Value2 = '  +11.51%  '
$ws.Range('E50').Value2 = '  +1.62%  '
$ws.Range('E51').Value2 = '  +0.40%  '
